$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "G2" = 0.7137273333333334
    "H2" = 2.141182
    "I2" = 0.3473329658610935
    "J2" = 0.3473329658610935
    "M2" = 34.30489066666667
    "N2" = 102.914672
    "O2" = 0.3410498283191481
    "P2" = 0.3410498283191481
    "Q2" = 24.48433813581156
    "R2" = 220.359043222304
    "S2" = 0.1184578483765065
    "T2" = 0.1184578483765065
    "G3" = 0.7137273333333334
    "H3" = 2.141182
    "I3" = 0.3473329658610935
    "J3" = 0.3473329658610935
    "O3" = 0.1375196956779303
    "P3" = 0.1375196956779303
    "Q3" = 9.872688533247112
    "R3" = 88.854196799224
    "S3" = 0.04776512376413054
    "T3" = 0.04776512376413054
    "G4" = 0.7137273333333334
    "H4" = 2.141182
    "I4" = 0.3473329658610935
    "J4" = 0.3473329658610935
    "M4" = 25.48814066666667
    "N4" = 76.46442200000001
    "O4" = 0.2533961143618365
    "P4" = 0.2533961143618365
    "Q4" = 18.19158266964489
    "R4" = 163.724244026804
    "S4" = 0.08801282393897349
    "T4" = 0.08801282393897351
    "G5" = 0.7137273333333334
    "H5" = 2.141182
    "I5" = 0.3473329658610935
    "J5" = 0.3473329658610935
    "M5" = 3.430646
    "N5" = 10.291938
    "O5" = 0.03410654302013726
    "P5" = 0.03410654302013727
    "Q5" = 2.448545821190667
    "R5" = 22.036912390716
    "S5" = 0.01184632674245325
    "T5" = 0.01184632674245326
    "G6" = 0.7137273333333334
    "H6" = 2.141182
    "I6" = 0.3473329658610935
    "J6" = 0.3473329658610935
    "M6" = 7.592334333333334
    "N6" = 22.777003
    "O6" = 0.07548090871605478
    "P6" = 0.0754809087160548
    "Q6" = 5.418856537505111
    "R6" = 48.76970883754601
    "S6" = 0.02621700789023777
    "T6" = 0.02621700789023778
    "G7" = 0.7137273333333334
    "H7" = 2.141182
    "I7" = 0.3473329658610935
    "J7" = 0.3473329658610935
    "M7" = 15.93756533333333
    "N7" = 47.812696
    "O7" = 0.158446909904893
    "P7" = 0.158446909904893
    "Q7" = 11.37507600518578
    "R7" = 102.375684046672
    "S7" = 0.05503383514879196
    "T7" = 0.05503383514879197
    "I8" = 0.2938237411507374
    "J8" = 0.2938237411507374
    "M8" = 34.30489066666667
    "N8" = 102.914672
    "O8" = 0.3410498283191481
    "P8" = 0.3410498283191481
    "Q8" = 20.71234388255823
    "R8" = 186.411094943024
    "S8" = 0.1002085364755488
    "T8" = 0.1002085364755488
    "I9" = 0.2938237411507374
    "J9" = 0.2938237411507374
    "O9" = 0.1375196956779303
    "P9" = 0.1375196956779303
    "S9" = 0.04040655146600038
    "T9" = 0.04040655146600038
    "I10" = 0.2938237411507374
    "J10" = 0.2938237411507374
    "M10" = 25.48814066666667
    "N10" = 76.46442200000001
    "O10" = 0.2533961143618365
    "P10" = 0.2533961143618365
    "Q10" = 15.38903416264156
    "R10" = 138.501307463774
    "S10" = 0.0744537943148549
    "T10" = 0.07445379431485491
    "I11" = 0.2938237411507374
    "J11" = 0.2938237411507374
    "M11" = 3.430646
    "N11" = 10.291938
    "O11" = 0.03410654302013726
    "P11" = 0.03410654302013727
    "Q11" = 2.071329140260667
    "R11" = 18.641962262346
    "S11" = 0.0100213120678953
    "T11" = 0.0100213120678953
    "I12" = 0.2938237411507374
    "J12" = 0.2938237411507374
    "M12" = 7.592334333333334
    "N12" = 22.777003
    "O12" = 0.07548090871605478
    "P12" = 0.0754809087160548
    "Q12" = 4.584041415883445
    "R12" = 41.256372742951
    "S12" = 0.02217808298440852
    "T12" = 0.02217808298440853
    "I13" = 0.2938237411507374
    "J13" = 0.2938237411507374
    "M13" = 15.93756533333333
    "N13" = 47.812696
    "O13" = 0.158446909904893
    "P13" = 0.158446909904893
    "Q13" = 9.62266100895911
    "R13" = 86.60394908063199
    "S13" = 0.04655546384202949
    "T13" = 0.04655546384202951
    "E14" = 1
    "F14" = 0.3333333333333333
    "G14" = 0.093007
    "H14" = 0.279021
    "I14" = 0.04526153847152096
    "J14" = 0.04526153847152095
    "M14" = 34.30489066666667
    "N14" = 102.914672
    "O14" = 0.3410498283191481
    "P14" = 0.3410498283191481
    "Q14" = 3.190594966234667
    "R14" = 28.715354696112
    "S14" = 0.01543643992517274
    "T14" = 0.01543643992517274
    "E15" = 1
    "F15" = 0.3333333333333333
    "G15" = 0.093007
    "H15" = 0.279021
    "I15" = 0.04526153847152096
    "J15" = 0.04526153847152095
    "O15" = 0.1375196956779303
    "P15" = 0.1375196956779303
    "Q15" = 1.286526520041333
    "R15" = 11.578738680372
    "S15" = 0.006224352996518496
    "T15" = 0.006224352996518497
    "E16" = 1
    "F16" = 0.3333333333333333
    "G16" = 0.093007
    "H16" = 0.279021
    "I16" = 0.04526153847152096
    "J16" = 0.04526153847152095
    "M16" = 25.48814066666667
    "N16" = 76.46442200000001
    "O16" = 0.2533961143618365
    "P16" = 0.2533961143618365
    "Q16" = 2.370575498984667
    "R16" = 21.335179490862
    "S16" = 0.01146909797872219
    "T16" = 0.01146909797872219
    "E17" = 1
    "F17" = 0.3333333333333333
    "G17" = 0.093007
    "H17" = 0.279021
    "I17" = 0.04526153847152096
    "J17" = 0.04526153847152095
    "M17" = 3.430646
    "N17" = 10.291938
    "O17" = 0.03410654302013726
    "P17" = 0.03410654302013727
    "Q17" = 0.319074092522
    "R17" = 2.871666832698
    "S17" = 0.001543714609036527
    "T17" = 0.001543714609036527
    "E18" = 1
    "F18" = 0.3333333333333333
    "G18" = 0.093007
    "H18" = 0.279021
    "I18" = 0.04526153847152096
    "J18" = 0.04526153847152095
    "M18" = 7.592334333333334
    "N18" = 22.777003
    "O18" = 0.07548090871605478
    "P18" = 0.0754809087160548
    "Q18" = 0.7061402393403334
    "R18" = 6.355262154063
    "S18" = 0.003416382053717075
    "T18" = 0.003416382053717075
    "E19" = 1
    "F19" = 0.3333333333333333
    "G19" = 0.093007
    "H19" = 0.279021
    "I19" = 0.04526153847152096
    "J19" = 0.04526153847152095
    "M19" = 15.93756533333333
    "N19" = 47.812696
    "O19" = 0.158446909904893
    "P19" = 0.158446909904893
    "Q19" = 1.482305138957333
    "R19" = 13.340746250616
    "S19" = 0.007171550908353929
    "T19" = 0.007171550908353929
    "G20" = 0.2640916666666667
    "H20" = 0.792275
    "I20" = 0.1285193064053396
    "J20" = 0.1285193064053396
    "M20" = 34.30489066666667
    "N20" = 102.914672
    "O20" = 0.3410498283191481
    "P20" = 0.3410498283191481
    "Q20" = 9.05963575097778
    "R20" = 81.5367217588
    "S20" = 0.04383148738523706
    "T20" = 0.04383148738523706
    "G21" = 0.2640916666666667
    "H21" = 0.792275
    "I21" = 0.1285193064053396
    "J21" = 0.1285193064053396
    "O21" = 0.1375196956779303
    "P21" = 0.1375196956779303
    "Q21" = 3.653068402255556
    "R21" = 32.8776156203
    "S21" = 0.01767393590560098
    "T21" = 0.01767393590560098
    "G22" = 0.2640916666666667
    "H22" = 0.792275
    "I22" = 0.1285193064053396
    "J22" = 0.1285193064053396
    "M22" = 25.48814066666667
    "N22" = 76.46442200000001
    "O22" = 0.2533961143618365
    "P22" = 0.2533961143618365
    "Q22" = 6.731205548894446
    "R22" = 60.58084994005001
    "S22" = 0.03256629286359133
    "T22" = 0.03256629286359134
    "G23" = 0.2640916666666667
    "H23" = 0.792275
    "I23" = 0.1285193064053396
    "J23" = 0.1285193064053396
    "M23" = 3.430646
    "N23" = 10.291938
    "O23" = 0.03410654302013726
    "P23" = 0.03410654302013727
    "Q23" = 0.9060050198833333
    "R23" = 8.15404517895
    "S23" = 0.004383349252831917
    "T23" = 0.004383349252831919
    "G24" = 0.2640916666666667
    "H24" = 0.792275
    "I24" = 0.1285193064053396
    "J24" = 0.1285193064053396
    "M24" = 7.592334333333334
    "N24" = 22.777003
    "O24" = 0.07548090871605478
    "P24" = 0.0754809087160548
    "Q24" = 2.005072227980556
    "R24" = 18.045650051825
    "S24" = 0.009700754035032113
    "T24" = 0.009700754035032115
    "G25" = 0.2640916666666667
    "H25" = 0.792275
    "I25" = 0.1285193064053396
    "J25" = 0.1285193064053396
    "M25" = 15.93756533333333
    "N25" = 47.812696
    "O25" = 0.158446909904893
    "P25" = 0.158446909904893
    "Q25" = 4.208978191488889
    "R25" = 37.88080372339999
    "S25" = 0.02036348696304618
    "T25" = 0.02036348696304618
    "E26" = 2
    "F26" = 0.6666666666666666
    "G26" = 0.2110656666666667
    "H26" = 0.633197
    "I26" = 0.1027143848511462
    "J26" = 0.1027143848511461
    "M26" = 34.30489066666667
    "N26" = 102.914672
    "O26" = 0.3410498283191481
    "P26" = 0.3410498283191481
    "Q26" = 7.240584618487112
    "R26" = 65.165261566384
    "S26" = 0.03503072331939031
    "T26" = 0.0350307233193903
    "E27" = 2
    "F27" = 0.6666666666666666
    "G27" = 0.2110656666666667
    "H27" = 0.633197
    "I27" = 0.1027143848511462
    "J27" = 0.1027143848511461
    "O27" = 0.1375196956779303
    "P27" = 0.1375196956779303
    "Q27" = 2.919582156578223
    "R27" = 26.276239409204
    "S27" = 0.01412525094647543
    "T27" = 0.01412525094647543
    "E28" = 2
    "F28" = 0.6666666666666666
    "G28" = 0.2110656666666667
    "H28" = 0.633197
    "I28" = 0.1027143848511462
    "J28" = 0.1027143848511461
    "M28" = 25.48814066666667
    "N28" = 76.46442200000001
    "O28" = 0.2533961143618365
    "P28" = 0.2533961143618365
    "Q28" = 5.379671401903779
    "R28" = 48.41704261713401
    "S28" = 0.02602742601034672
    "T28" = 0.02602742601034672
    "E29" = 2
    "F29" = 0.6666666666666666
    "G29" = 0.2110656666666667
    "H29" = 0.633197
    "I29" = 0.1027143848511462
    "J29" = 0.1027143848511461
    "M29" = 3.430646
    "N29" = 10.291938
    "O29" = 0.03410654302013726
    "P29" = 0.03410654302013727
    "Q29" = 0.7240915850873334
    "R29" = 6.516824265786
    "S29" = 0.003503232585712551
    "T29" = 0.003503232585712552
    "E30" = 2
    "F30" = 0.6666666666666666
    "G30" = 0.2110656666666667
    "H30" = 0.633197
    "I30" = 0.1027143848511462
    "J30" = 0.1027143848511461
    "M30" = 7.592334333333334
    "N30" = 22.777003
    "O30" = 0.07548090871605478
    "P30" = 0.0754809087160548
    "Q30" = 1.602481107621222
    "R30" = 14.422329968591
    "S30" = 0.007752975106775084
    "T30" = 0.007752975106775084
    "E31" = 2
    "F31" = 0.6666666666666666
    "G31" = 0.2110656666666667
    "H31" = 0.633197
    "I31" = 0.1027143848511462
    "J31" = 0.1027143848511461
    "M31" = 15.93756533333333
    "N31" = 47.812696
    "O31" = 0.158446909904893
    "P31" = 0.158446909904893
    "Q31" = 3.363872852123555
    "R31" = 30.274855669112
    "S31" = 0.01627477688244606
    "T31" = 0.01627477688244606
    "G32" = 0.1692153333333334
    "H32" = 0.507646
    "I32" = 0.08234806326016224
    "J32" = 0.08234806326016222
    "M32" = 34.30489066666667
    "N32" = 102.914672
    "O32" = 0.3410498283191481
    "P32" = 0.3410498283191481
    "Q32" = 5.804913509123558
    "R32" = 52.24422158211201
    "S32" = 0.02808479283729268
    "T32" = 0.02808479283729267
    "G33" = 0.1692153333333334
    "H33" = 0.507646
    "I33" = 0.08234806326016224
    "J33" = 0.08234806326016222
    "O33" = 0.1375196956779303
    "P33" = 0.1375196956779303
    "Q33" = 2.340684184319112
    "R33" = 21.066157658872
    "S33" = 0.01132448059920446
    "T33" = 0.01132448059920446
    "G34" = 0.1692153333333334
    "H34" = 0.507646
    "I34" = 0.08234806326016224
    "J34" = 0.08234806326016222
    "M34" = 25.48814066666667
    "N34" = 76.46442200000001
    "O34" = 0.2533961143618365
    "P34" = 0.2533961143618365
    "Q34" = 4.31298421895689
    "R34" = 38.81685797061201
    "S34" = 0.02086667925534781
    "T34" = 0.02086667925534781
    "G35" = 0.1692153333333334
    "H35" = 0.507646
    "I35" = 0.08234806326016224
    "J35" = 0.08234806326016222
    "M35" = 3.430646
    "N35" = 10.291938
    "O35" = 0.03410654302013726
    "P35" = 0.03410654302013727
    "Q35" = 0.5805179064386667
    "R35" = 5.224661157948001
    "S35" = 0.002808607762207707
    "T35" = 0.002808607762207708
    "G36" = 0.1692153333333334
    "H36" = 0.507646
    "I36" = 0.08234806326016224
    "J36" = 0.08234806326016222
    "M36" = 7.592334333333334
    "N36" = 22.777003
    "O36" = 0.07548090871605478
    "P36" = 0.0754809087160548
    "Q36" = 1.284739384993111
    "R36" = 11.562654464938
    "S36" = 0.00621570664588421
    "T36" = 0.00621570664588421
    "G37" = 0.1692153333333334
    "H37" = 0.507646
    "I37" = 0.08234806326016224
    "J37" = 0.08234806326016222
    "M37" = 15.93756533333333
    "N37" = 47.812696
    "O37" = 0.158446909904893
    "P37" = 0.158446909904893
    "Q37" = 2.696880430401778
    "R37" = 24.271923873616
    "S37" = 0.01304779616022535
    "T37" = 0.01304779616022535
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
